# Update stats for 2025-08 (row 21 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6236
$ws.Range("D21").Value = 5592802
$ws.Range("E21").Value = 896.8572803078897
$ws.Range("F21").Value = 8.245096337441415
$ws.Range("H21").Value = 27.66121308239358
